$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.23860000000001
$ws.Range("E3").Value = 16.0294
$ws.Range("A12").Value = -21.53540000000001
$ws.Range("C14").Value = -13.1498
$ws.Range("C26").Value = -12.3412
$ws.Range("E30").Value = 15.9907
$ws.Range("C31").Value = -12.85740000000001
$ws.Range("A32").Value = -21.34199999999998
$ws.Range("C35").Value = -12.85970000000001
$ws.Range("A36").Value = -20.0946
$ws.Range("C37").Value = -14.073
$ws.Range("A38").Value = -19.50829999999999
$ws.Range("E44").Value = 16.8749
$ws.Range("C45").Value = -13.68019999999999
$ws.Range("A46").Value = -21.61169999999999
$ws.Range("A54").Value = -21.87009999999998
$ws.Range("A55").Value = -22.34980000000001
$ws.Range("C57").Value = -14.33729999999998
$ws.Range("E58").Value = 16.35840000000001
$ws.Range("A67").Value = -21.50379999999997
$ws.Range("A69").Value = -21.63499999999998
$ws.Range("A72").Value = -21.68469999999999
$ws.Range("E84").Value = 16.65720000000001
$ws.Range("E89").Value = 17.29640000000002
$ws.Range("A91").Value = -21.44300000000001
$ws.Range("E91").Value = 17.97820000000002
$ws.Range("E92").Value = 18.07810000000001
$ws.Range("A99").Value = -20.23089999999999
$ws.Range("C100").Value = -12.7328
$ws.Range("C102").Value = -12.37479999999999
$ws.Range("E102").Value = 16.74239999999999
